# Commit: "Committed on 08th March 2017"
# Adds two new rows of the "day 1" log (rows 99 & 100), fills in the
# previously-empty Time/Errors/Error&Solution cells of row 98, and moves
# the sheet's active-cell selection down to A100.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day 1")

# --- Row 98: finish out the Time Taken / Errors Y-N / Error&Solution cells
$ws.Range("F98").Value = "10 minutes"
$ws.Range("G98").Value = "N"
$ws.Range("H98").Value = "NA"

# --- Row 99: new entry - "Creation of RegisterHandler"
$ws.Range("C99").Value = "Creation of RegisterHandler"
$ws.Range("B99").Value = "8th Mar,2017"
$ws.Range("D99").Value = "NA"
$ws.Range("E99").Value = "NA"
$ws.Range("F99").Value = "20 minutes"
$ws.Range("G99").Value = "N"
$ws.Range("H99").Value = "NA"
$ws.Rows.Item(99).RowHeight = 28.8

# --- Row 100: new entry - "Creation of Excepiton.jsp page"
$ws.Range("B100").Value = "8th Mar,2017"
$ws.Range("C100").Value = "Creation of Excepiton.jsp page"
$ws.Range("D100").Value = "NA"
$ws.Range("E100").Value = "NA"
$ws.Range("F100").Value = "10 minutes"
$ws.Range("G100").Value = "N"
$ws.Range("H100").Value = "NA"
$ws.Rows.Item(100).RowHeight = 28.8

# --- Move the active selection to A100 (matches the author's final click)
[void]$ws.Range("A100").Select()
